$d = $word.ActiveDocument

# 1. Insert "FrankAbba, " before "GalloSalvato" with correct formatting
$findRng = $d.Content
$findRng.Find.Execute("GalloSalvato", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos = $findRng.Start

# Insert ", " first so it inherits plain formatting from the preceding space run
$commaRng = $d.Range($insPos, $insPos)
$commaRng.Text = ", "

# Insert "FrankAbba" before the comma and make it bold+italic
$nameRng = $d.Range($insPos, $insPos)
$nameRng.Text = "FrankAbba"
$nameRng.Bold = 1
$nameRng.Italic = 1

# 2. Move the _GoBack bookmark into "opera è generato " splitting it into
#    "opera è gen" | bookmark | "erato "  (done after the text insertion above,
#    matching the order of edits actually performed by the author)
$bmFindRng = $d.Content
$bmFindRng.Find.Execute("opera è gen", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $bmFindRng.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
